$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update existing rows 59-65 (values are re-ordered/changed per source diff)
# ---------------------------------------------------------------------------

# Row 59
$ws.Range("D59").Value = 45015
$ws.Range("L59").Value = "Especial"
$ws.Range("M59").Value = 180
$ws.Range("N59").Value = 10000
$ws.Range("O59").Value = 10000
$ws.Range("P59").Value = 10000
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 556

# Row 60
$ws.Range("D60").Value = 45015
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 220
$ws.Range("N60").Value = 8000
$ws.Range("O60").Value = 8000
$ws.Range("P60").Value = 8000
$ws.Range("S60").Value = 444

# Row 61
$ws.Range("D61").Value = 44291
$ws.Range("M61").Value = 300
$ws.Range("N61").Value = 9000
$ws.Range("O61").Value = 9000
$ws.Range("P61").Value = 9000
$ws.Range("R61").Value = "Provincia de Curicó"
$ws.Range("S61").Value = 500

# Row 62
$ws.Range("D62").Value = 44998
$ws.Range("M62").Value = 200
$ws.Range("N62").Value = 12000
$ws.Range("O62").Value = 12000
$ws.Range("P62").Value = 12000
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 667

# Row 63
$ws.Range("D63").Value = 44998
$ws.Range("M63").Value = 250
$ws.Range("N63").Value = 10000
$ws.Range("O63").Value = 10000
$ws.Range("P63").Value = 10000
$ws.Range("R63").Value = "Región de O'Higgins"
$ws.Range("S63").Value = 556

# Row 64
$ws.Range("D64").Value = 45013

# Row 65
$ws.Range("D65").Value = 45013
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 8000
$ws.Range("O65").Value = 8000
$ws.Range("P65").Value = 8000
$ws.Range("Q65").Value = "$/caja 18 kilos granel"
$ws.Range("R65").Value = "Provincia de Curicó"
$ws.Range("S65").Value = 444
$ws.Range("T65").Value = 18

# ---------------------------------------------------------------------------
# Insert two brand-new rows (66 and 67), pushing the former row 66 to row 68
# ---------------------------------------------------------------------------

$ws.Rows.Item(66).Insert()

$ws.Range("A66").Value = 5
$ws.Range("B66").Value = "Macroferia Regional de Talca"
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = 44306
$ws.Range("E66").Value = 7
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100104
$ws.Range("H66").Value = "Frutos de pepita"
$ws.Range("I66").Value = 100104003
$ws.Range("J66").Value = "Membrillo"
$ws.Range("K66").Value = "Champion"
$ws.Range("L66").Value = "Especial"
$ws.Range("M66").Value = 230
$ws.Range("N66").Value = 10000
$ws.Range("O66").Value = 10000
$ws.Range("P66").Value = 10000
$ws.Range("Q66").Value = "$/caja 18 kilos granel"
$ws.Range("R66").Value = "Provincia de Curicó"
$ws.Range("S66").Value = 556
$ws.Range("T66").Value = 18

$ws.Rows.Item(67).Insert()

$ws.Range("A67").Value = 5
$ws.Range("B67").Value = "Macroferia Regional de Talca"
$ws.Range("C67").Value = "Maule"
$ws.Range("D67").Value = 44342
$ws.Range("E67").Value = 7
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100104
$ws.Range("H67").Value = "Frutos de pepita"
$ws.Range("I67").Value = 100104003
$ws.Range("J67").Value = "Membrillo"
$ws.Range("K67").Value = "Champion"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 300
$ws.Range("N67").Value = 9000
$ws.Range("O67").Value = 9000
$ws.Range("P67").Value = 9000
$ws.Range("Q67").Value = "$/caja 15 kilos empedrada"
$ws.Range("R67").Value = "Región de O'Higgins"
$ws.Range("S67").Value = 600
$ws.Range("T67").Value = 15

# Ensure column D keeps the date number-format/style used throughout the sheet
$ws.Range("D66").NumberFormat = $ws.Range("D65").NumberFormat
$ws.Range("D67").NumberFormat = $ws.Range("D65").NumberFormat
